# Updated symbol list on Sat Jan  7 12:43:28 UTC 2023 with GitHub Actions
#
# Refresh cryptocurrency price / 1h-volume snapshot values on Sheet1.
# Each target cell originally holds its value as literal text (e.g. "260.59",
# "1.74%"), so we force the Text number format ("@") before assigning the new
# string value, then restore the "Normal" style so no stray style index is
# left behind (matches the source workbook's un-styled D/E data cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "260.21"
Set-TextValue "E2" "1.85%"
Set-TextValue "D3" "27.28"
Set-TextValue "E3" "2.56%"
Set-TextValue "D4" "4.678"
Set-TextValue "E4" "0.71%"
Set-TextValue "D5" "0.06100"
Set-TextValue "E5" "2.72%"
Set-TextValue "D6" "6.662"
Set-TextValue "E6" "0.64%"
Set-TextValue "D7" "0.8503"
Set-TextValue "E8" "0.71%"
Set-TextValue "D9" "0.1397"
Set-TextValue "E9" "1.16%"
Set-TextValue "D10" "0.04856"
Set-TextValue "E10" "19.08%"
Set-TextValue "D11" "0.07091"
Set-TextValue "E11" "1.37%"
Set-TextValue "D12" "0.03072"
Set-TextValue "E12" "0.50%"
Set-TextValue "D13" "0.09051"
Set-TextValue "E13" "-0.51%"
Set-TextValue "D14" "0.001528"
Set-TextValue "E14" "0.44%"
Set-TextValue "D15" "0.0006079"
Set-TextValue "E15" "-94.06%"
Set-TextValue "D16" "0.006162"
Set-TextValue "E16" "1.95%"
Set-TextValue "D17" "3.451"
Set-TextValue "E17" "-0.44%"
Set-TextValue "D18" "3.152"
Set-TextValue "E18" "0.34%"
Set-TextValue "E19" "-0.62%"
Set-TextValue "E20" "2.97%"
Set-TextValue "E21" "1.56%"
Set-TextValue "D22" "4.091"
Set-TextValue "E22" "6.13%"
Set-TextValue "D23" "0.04243"
Set-TextValue "E23" "0.76%"
Set-TextValue "D24" "0.001222"
Set-TextValue "E24" "0.51%"
Set-TextValue "D25" "0.003800"
Set-TextValue "E25" "-19.38%"
Set-TextValue "E26" "0.06%"
Set-TextValue "E27" "3.42%"
Set-TextValue "D40" "0.03856"
Set-TextValue "E40" "2.36%"
Set-TextValue "E41" "1.59%"
Set-TextValue "D42" "0.004080"
Set-TextValue "E42" "-34.61%"
Set-TextValue "D43" "0.01631"
Set-TextValue "E43" "13.58%"
Set-TextValue "D44" "0.002218"
Set-TextValue "E44" "0.82%"
Set-TextValue "D45" "0.00005155"
Set-TextValue "E45" "1.30%"
Set-TextValue "E46" "0.06%"
Set-TextValue "D47" "0.1371"
Set-TextValue "E47" "-43.08%"
Set-TextValue "E48" "8.92%"
Set-TextValue "E49" "0.06%"
Set-TextValue "E50" "0.06%"
